$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.507.42'
$ws.Cells.Item(2, 5).Value = '  -0.45%  '
$ws.Cells.Item(3, 4).Value = '2.518.61'
$ws.Cells.Item(3, 5).Value = '  -0.99%  '
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  -0.07%  '
$ws.Cells.Item(5, 4).Value = '''311.65'
$ws.Cells.Item(5, 5).Value = '  -0.53%  '
$ws.Cells.Item(6, 4).Value = '''98.80'
$ws.Cells.Item(6, 5).Value = '  -2.40%  '
$ws.Cells.Item(7, 5).Value = '  -1.26%  '
$ws.Cells.Item(8, 5).Value = '  -0.03%  '
$ws.Cells.Item(9, 5).Value = '  -2.93%  '
$ws.Cells.Item(10, 4).Value = '''35.20'
$ws.Cells.Item(10, 5).Value = '  -2.66%  '
$ws.Cells.Item(11, 4).Value = '''0.0801'
$ws.Cells.Item(11, 5).Value = '  -0.70%  '
$ws.Cells.Item(12, 5).Value = '  +0.21%  '
$ws.Cells.Item(13, 4).Value = '''7.21'
$ws.Cells.Item(13, 5).Value = '  -2.60%  '
$ws.Cells.Item(14, 4).Value = '2.900.00'
$ws.Cells.Item(14, 5).Value = '  -1.30%  '
$ws.Cells.Item(15, 4).Value = '2.597.38'
$ws.Cells.Item(15, 5).Value = '  +3.82%  '
$ws.Cells.Item(16, 4).Value = '''15.28'
$ws.Cells.Item(16, 5).Value = '  -3.44%  '
$ws.Cells.Item(17, 4).Value = '''0.808'
$ws.Cells.Item(17, 5).Value = '  -3.60%  '
$ws.Cells.Item(18, 4).Value = '42.486.30'
$ws.Cells.Item(18, 5).Value = '  -0.58%  '
$ws.Cells.Item(19, 5).Value = '  -2.57%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0945'
$ws.Cells.Item(20, 5).Value = '  -1.04%  '
$ws.Cells.Item(21, 4).Value = '''11.93'
$ws.Cells.Item(21, 5).Value = '  -3.43%  '
$ws.Cells.Item(22, 4).Value = '''68.89'
$ws.Cells.Item(22, 5).Value = '  -0.32%  '
$ws.Cells.Item(23, 4).Value = '''240.64'
$ws.Cells.Item(23, 5).Value = '  -2.22%  '
$ws.Cells.Item(24, 4).Value = '''2.84'
$ws.Cells.Item(24, 5).Value = '  -3.38%  '
$ws.Cells.Item(25, 5).Value = '  -3.58%  '
$ws.Cells.Item(26, 5).Value = '  +0.12%  '
$ws.Cells.Item(27, 5).Value = '  -4.87%  '
$ws.Cells.Item(28, 5).Value = '  -4.11%  '
$ws.Cells.Item(29, 5).Value = '  -1.09%  '
$ws.Cells.Item(30, 4).Value = '''38.40'
$ws.Cells.Item(30, 5).Value = '  -6.14%  '
$ws.Cells.Item(31, 5).Value = '  +2.24%  '
$ws.Cells.Item(32, 4).Value = '''156.35'
$ws.Cells.Item(32, 5).Value = '  -0.21%  '
$ws.Cells.Item(33, 4).Value = '''2.80'
$ws.Cells.Item(33, 5).Value = '  +5.93%  '
$ws.Cells.Item(34, 5).Value = '  +1.22%  '
$ws.Cells.Item(35, 4).Value = '''0.0784'
$ws.Cells.Item(35, 5).Value = '  -2.43%  '
$ws.Cells.Item(36, 5).Value = '  -4.03%  '
$ws.Cells.Item(37, 5).Value = '  -5.98%  '
$ws.Cells.Item(38, 4).Value = '''17.44'
$ws.Cells.Item(38, 5).Value = '  -4.74%  '
$ws.Cells.Item(39, 5).Value = '  -3.71%  '
$ws.Cells.Item(40, 5).Value = '  -0.85%  '
$ws.Cells.Item(41, 4).Value = '''4.17'
$ws.Cells.Item(41, 5).Value = '  -0.41%  '
$ws.Cells.Item(42, 4).Value = '''21.65'
$ws.Cells.Item(42, 5).Value = '  -3.04%  '
$ws.Cells.Item(43, 5).Value = '  -0.08%  '
$ws.Cells.Item(44, 5).Value = '  -1.01%  '
$ws.Cells.Item(45, 5).Value = '  -1.42%  '
$ws.Cells.Item(46, 4).Value = '1.996.64'
$ws.Cells.Item(46, 5).Value = '  +0.77%  '
$ws.Cells.Item(47, 5).Value = '  +0.65%  '
$ws.Cells.Item(48, 4).Value = '2.757.38'
$ws.Cells.Item(48, 5).Value = '  -1.29%  '
$ws.Cells.Item(49, 5).Value = '  -2.20%  '
$ws.Cells.Item(50, 4).Value = '''78.68'
$ws.Cells.Item(50, 5).Value = '  -3.32%  '
$ws.Cells.Item(51, 4).Value = '''100.21'
$ws.Cells.Item(51, 5).Value = '  -1.51%  '
